# Add the new Beneficiary Location roll-up fields (Latitude, Longitude,
# Geographic Code, Geographic Code Type) to the "Activity" summary sheet.
# They are inserted right after the existing "Beneficiary Location:Country
# Code" column, shifting every later column (Funding Org, Grant Programme,
# Related Activity, Last modified, Data Source, ...) four places to the
# right.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Activity")

# "Beneficiary Location:Country Code" currently lives in column W (23rd
# column); insert four new columns immediately after it, shifting the
# remaining header cells (and the used range) to the right.
$ws.Range("X1:AA1").Insert(-4161)

$ws.Range("X1").Value = "Beneficiary Location:Latitude"
$ws.Range("Y1").Value = "Beneficiary Location:Longitude"
$ws.Range("Z1").Value = "Beneficiary Location:Geographic Code"
$ws.Range("AA1").Value = "Beneficiary Location:Geographic Code Type"
